# Commit: "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" (column E) / "Valor Mora" (column F) rows for the
# worker's account statement (rows 16-22) get their period order reversed:
# the newest period (2209) now sits on top (row 16) and the oldest (2203)
# sits at the bottom (row 22); the reduced/partial payment amount (34666)
# moves from the old period (2209, was row 22) to the new period (2209,
# now row 16).
#
# We read the current column E/F values first, then write them back in
# reversed order so each row keeps its own formatting/style untouched and
# only the underlying data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 22

# Capture existing "Periodo Mora" (E) values for rows 16..22
$periodo = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodo[$r] = $ws.Cells.Item($r, 5).Value2
}

# Capture existing "Valor Mora" (F) values for rows 16..22
$valorMora = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $valorMora[$r] = $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse row order (row 16 <-> row 22, 17 <-> 21, 18 <-> 20, 19 stays)
$n = $lastRow - $firstRow
for ($i = 0; $i -le $n; $i++) {
    $srcRow = $firstRow + $i
    $dstRow = $lastRow - $i

    $ws.Cells.Item($dstRow, 5).Value = $periodo[$srcRow]
    $ws.Cells.Item($dstRow, 6).Value = $valorMora[$srcRow]
}
